$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.368.80"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.07%  "

$ws.Range("D3").Value = "'1.938.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -3.05%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.17%  "

$ws.Range("D5").Value = "'251.14"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.13%  "

$ws.Range("D6").Value = "'0.7086"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.79%  "

$ws.Range("E7").Value = "  -0.09%  "

$ws.Range("D8").Value = "'0.3301"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.67%  "

$ws.Range("E9").Value = "  -1.02%  "

$ws.Range("D10").Value = "'0.07306"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.34%  "

$ws.Range("D11").Value = "'0.8057"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.81%  "

$ws.Range("D12").Value = "'0.08083"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.92%  "

$ws.Range("D13").Value = "'1.935.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.27%  "

$ws.Range("D14").Value = "'5.486"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.78%  "

$ws.Range("D15").Value = "'94.52"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -6.35%  "

$ws.Range("D16").Value = "'15.05"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -3.97%  "

$ws.Range("D17").Value = "'30.350.75"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.13%  "

$ws.Range("D18").Value = "'253.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -5.42%  "

$ws.Range("D19").Value = "'0.000008195"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.71%  "

$ws.Range("D20").Value = "'5.816"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.88%  "

$ws.Range("D21").Value = "'2.190.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.24%  "

$ws.Range("D22").Value = "'1.003"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.06%  "

$ws.Range("D23").Value = "'1.004"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.04%  "

$ws.Range("D24").Value = "'7.009"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.61%  "

$ws.Range("D25").Value = "'9.729"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.65%  "

$ws.Range("D26").Value = "'164.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.11%  "

$ws.Range("D27").Value = "'2.351"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.89%  "

$ws.Range("D28").Value = "'19.31"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.58%  "

$ws.Range("D29").Value = "'0.1290"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.26%  "

$ws.Range("E30").Value = "  -2.04%  "

$ws.Range("E31").Value = "  -3.44%  "

$ws.Range("D32").Value = "'4.424"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -5.24%  "

$ws.Range("D33").Value = "'4.165"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -6.66%  "

$ws.Range("D34").Value = "'0.05183"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.37%  "

$ws.Range("D35").Value = "'1.265"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -2.47%  "

$ws.Range("D36").Value = "'0.7473"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -5.58%  "

$ws.Range("D37").Value = "'2.751"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.44%  "

$ws.Range("D38").Value = "'0.01968"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.60%  "

$ws.Range("E39").Value = "  -3.50%  "

$ws.Range("D40").Value = "'78.83"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.62%  "

$ws.Range("D41").Value = "'6.428"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.84%  "

$ws.Range("D42").Value = "'0.4532"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.45%  "

$ws.Range("D43").Value = "'2.015"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.83%  "

$ws.Range("D44").Value = "'0.8472"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.59%  "

$ws.Range("D45").Value = "'1.003"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.06%  "

$ws.Range("D46").Value = "'101.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.79%  "

$ws.Range("D47").Value = "'9.772"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.24%  "

$ws.Range("D48").Value = "'7.440"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.57%  "

$ws.Range("D49").Value = "'36.63"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.75%  "

$ws.Range("D50").Value = "'0.4181"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.56%  "

$ws.Range("D51").Value = "'0.06038"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.54%  "
